$wb = $excel.ActiveWorkbook

# --- Update localization status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: zh-cn / de-de status columns (E, F) for both data rows
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the now shorter "Status" columns to fit the new text ---
# (ColumnWidth is quantized internally to 1/6-character steps, so 12.5 is the
# closest achievable setting to the recorded target width.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
